$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timelog")

# Update the hours worked for week row 14 (Monday=C, Tuesday=D, Wednesday=E)
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 2.5

# Move the active selection to F14
$ws.Range("F14").Select()
